# Update cryptocurrency price/volume data on Sheet1 (rows 2-51, columns D "Price" and E "Volume(1h)")
# matching the Thu Jul 13 07:44:03 UTC 2023 GitHub Actions data refresh.
#
# Column D holds price figures rendered as plain text (e.g. thousand-separated
# "30.391.78" or decimal "0.2872"). Purely numeric-looking values are forced to
# Text format before the assignment so Excel keeps the exact original digits
# (including trailing zeros / long decimals) instead of re-parsing them as
# floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.391.78'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = '1.871.97'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.56'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("E7").Value = '  -1.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2872'
$ws.Range("E8").Value = '  -2.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06492'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("E10").Value = '  -1.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '100.15'
$ws.Range("E11").Value = '  +2.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07795'
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").Value = '1.871.60'
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7264'
$ws.Range("E14").Value = '  -1.76%  '
$ws.Range("E15").Value = '  -1.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '285.03'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").Value = '30.379.25'
$ws.Range("E17").Value = '  -1.56%  '
$ws.Range("E18").Value = '  -1.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9998'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007496'
$ws.Range("E20").Value = '  -1.29%  '
$ws.Range("D21").Value = '2.115.64'
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.345'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.314'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '163.42'
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("E26").Value = '  -2.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.96'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.895'
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.09667'
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.321'
$ws.Range("E30").Value = '  -1.84%  '
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.222'
$ws.Range("E32").Value = '  -2.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.142'
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04813'
$ws.Range("E34").Value = '  -1.79%  '
$ws.Range("E35").Value = '  -0.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6885'
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01900'
$ws.Range("E38").Value = '  -0.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.842'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '76.28'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.281'
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.959'
$ws.Range("E42").Value = '  -3.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4218'
$ws.Range("E43").Value = '  -1.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8217'
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.10'
$ws.Range("E46").Value = '  -0.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.783'
$ws.Range("E47").Value = '  +3.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.010'
$ws.Range("E48").Value = '  -1.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.09'
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05767'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '888.48'
$ws.Range("E51").Value = '  -4.34%  '
